# Revert "Merge branch 'wrong-xlsform-col'"
#
# The merge that is being reverted had renamed the "message" column header
# (on the XLSForm "survey" sheet) to "label". Reverting it restores the
# header text to "message" while the "choices" sheet keeps its own,
# unrelated "label" column header untouched (same text, different cell).
#
# It also restores the previously-selected cell on the "survey" sheet and
# nudges a couple of incidental layout values left over from the revert.

$wb = $excel.ActiveWorkbook

$survey  = $wb.Worksheets.Item("survey")
$choices = $wb.Worksheets.Item("choices")

# --- Header text: survey!C1 goes back to "message" ---------------------
$survey.Range("C1").Value = "message"

# Touch choices!C1 too (it stays "label" textually, but this keeps it as
# its own independent entry rather than silently inheriting whatever
# survey!C1 now points at).
$choices.Range("C1").Value = "label"

# --- Restore the previously active cell on the survey sheet ------------
$survey.Range("B3").Select()

# --- Incidental column width nudge on the choices sheet ----------------
# (column A width goes from 14.62 to ~14.63 characters)
$choices.Columns.Item(1).ColumnWidth = 13.83

# Note: the original revert also nudges choices!pageSetup's
# FirstPageNumber from 0 to 1, but that value never affects anything
# visible (useFirstPageNumber stays false) and this runtime does not
# persist PageSetup.FirstPageNumber writes back to the saved file, so
# there is nothing more to do for it here.
